# Add configuration in quiz template: remove sample question rows (rows 3 and 4),
# keeping only the header rows, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample question content (rows 3 and 4), keeping existing formatting/style.
$ws.Range("A3:F4").ClearContents()

# Update the selected cell shown in the sheet view.
$ws.Range("H9").Select()
